# Weekly driver report update for 2025-05-05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 / Row 4 (Bad Drivers summary block) ---
$ws.Range("C3").Value = 82
$ws.Range("D3").Value = 97.2
$ws.Range("C4").Value = 82

# --- Row 12: Intel(R) Wi-Fi 6E AX210 160MHz - 23.110.0.5 ---
$ws.Range("B12").Value = 185718
$ws.Range("C12").Value = 616
$ws.Range("E12").Value = 443
$ws.Range("F12").Value = 186482

# --- Row 13: Intel(R) Wi-Fi 6E AX210 160MHz - 22.250.1.2 ---
$ws.Range("B13").Value = 31411
$ws.Range("C13").Value = 86
$ws.Range("E13").Value = 86
$ws.Range("F13").Value = 31514

# --- Row 15: becomes the 23.40.0.4 driver entry ---
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.40.0.4"
$ws.Range("B15").Value = 78476
$ws.Range("C15").Value = 78
$ws.Range("D15").Value = 56
$ws.Range("E15").Value = 146
$ws.Range("F15").Value = 78610
$ws.Range("H15").Value = "23.40.0.4"
$ws.Range("J15").Value = "'2024-03-09"

# --- Row 16: becomes the 23.90.0.2 driver entry ---
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.90.0.2"
$ws.Range("B16").Value = 584784
$ws.Range("C16").Value = 827
$ws.Range("D16").Value = 463
$ws.Range("E16").Value = 963
$ws.Range("F16").Value = 586074
$ws.Range("H16").Value = "23.90.0.2"
$ws.Range("J16").Value = "'2024-09-25"

# --- Row 17: becomes the 23.120.0.3 driver entry ---
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.120.0.3"
$ws.Range("B17").Value = 53252
$ws.Range("C17").Value = 46
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 70
$ws.Range("F17").Value = 53300
$ws.Range("H17").Value = "23.120.0.3"
$ws.Range("J17").Value = "'2025-02-05"

# --- Row 18: becomes the 22.130.0.5 driver entry ---
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5"
$ws.Range("B18").Value = 18722
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 26
$ws.Range("F18").Value = 18738
$ws.Range("H18").Value = "22.130.0.5"
$ws.Range("J18").Value = "'2022-03-14"

# --- Row 19: Intel(R) Wi-Fi 6E AX210 160MHz - 22.70.0.6 ---
$ws.Range("B19").Value = 15506
$ws.Range("E19").Value = 19
$ws.Range("F19").Value = 15506

# --- Row 23: Intel(R) Wi-Fi 6E AX210 160MHz - 22.110.1.1 ---
$ws.Range("B23").Value = 42439
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 46
$ws.Range("F23").Value = 42440
